$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new column F for "ISBN13" (shifts old F..I -> G..J) ---
$ws.Columns("F:F").Insert()

# --- Header cell F1: reuse the existing bold header font (same as E1),
#     then give it an integer number format, and set its text ---
$scratchHeader = $ws.Range("Z99")
$ws.Range("E1").Copy()
$scratchHeader.PasteSpecial(-4122)
$scratchHeader.NumberFormat = "0"
$scratchHeader.Copy()
$ws.Range("F1").PasteSpecial(-4122)
$scratchHeader.Clear()
$ws.Range("F1").Value = "ISBN13"

# --- Body cells F2:F11: Arial / #333333, integer number format ---
$scratchBody = $ws.Range("Z100")
$scratchBody.Font.Name = "Arial"
$scratchBody.Font.Color = 3355443
$scratchBody.NumberFormat = "0"
$scratchBody.Copy()
$ws.Range("F2:F11").PasteSpecial(-4122)
$scratchBody.Clear()

$ws.Range("F2").Value = 9780205309023
$ws.Range("F3").Value = 9781137585042
$ws.Range("F4").Value = 9781743214404
$ws.Range("F5").Value = 9780399162718
$ws.Range("F6").Value = 9780521189064
$ws.Range("F7").Value = 9780071453875
$ws.Range("F8").Value = 9780060891541
$ws.Range("F9").Value = 9781472244444
$ws.Range("F10").Value = 9781447261131
$ws.Range("F11").Value = 9781509814756

# --- Pre-existing Arial/#333333 cell (B10, title with wrap text) loses its bold ---
$ws.Range("B10").Font.Bold = $false
# Row 10 no longer needs its tall custom height now that content changed
$ws.Rows("10:10").AutoFit()

# --- Column widths (approximate character widths closest to target) ---
$ws.Columns("B:B").ColumnWidth = 111.66666666666667
$ws.Columns("C:C").ColumnWidth = 19.5
$ws.Columns("D:D").ColumnWidth = 17.833333333333332
$ws.Columns("E:E").ColumnWidth = 8.666666666666666
$ws.Columns("F:F").ColumnWidth = 15.166666666666666
$ws.Columns("G:G").ColumnWidth = 12.0
$ws.Columns("H:H").ColumnWidth = 15.666666666666666
$ws.Columns("I:I").ColumnWidth = 22.5
$ws.Columns("J:J").ColumnWidth = 33.5

# --- Selection / view state ---
$ws.Range("G15").Select()
